$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "thermal men leggings"
$ws.Range("A2").Value = "compression pants with knee pads"
$ws.Range("A3").Value = "under armor compression pants men"
$ws.Range("A4").Value = "compression pants with pads"
$ws.Range("A5").Value = "basketball warm up pants men"
$ws.Range("A6").Value = "basketball compression leggings"
$ws.Range("A7").Value = "bjj compression pants"
$ws.Range("A8").Value = "padded leggings basketball"
$ws.Range("A9").Value = "football tights"
$ws.Range("A10").Value = "knee pads for volleyball youth"
$ws.Range("A11").Value = "recovery compression pants men"
$ws.Range("A12").Value = "hockey compression leggings"
$ws.Range("A13").Value = "knee compressions"
$ws.Range("A14").Value = "men sports tights"
$ws.Range("A15").Value = "wrestling knee pads men"
$ws.Range("A16").Value = "mens running leggings"
$ws.Range("A17").Value = "thick knee pad"
$ws.Range("A18").Value = "football pads pants"
$ws.Range("A19").Value = "boys athletic pants"
$ws.Range("A20").Value = "running compression tights"
$ws.Range("A21").Value = "volleyball gear men"
$ws.Range("A22").Value = "girls capri leggings"
$ws.Range("A23").Value = "black capri leggings"
$ws.Range("A24").Value = "coolomg basketball knee pads"
$ws.Range("A25").Value = "tight with knee pads"
$ws.Range("A26").Value = "tights mens nike"
$ws.Range("A27").Value = "muscle compression pants for men"
$ws.Range("A28").Value = "hayabusa compression pants men"
$ws.Range("A29").Value = "girls athletic leggings"
$ws.Range("A30").Value = "black capri leggings for women"
$ws.Range("A31").Value = "asics youth knee pads"
$ws.Range("A32").Value = "men's tights sports"
$ws.Range("A33").Value = "tights with knee pads"
$ws.Range("A34").Value = "lavento compression pants"
$ws.Range("A35").Value = "tsla compression pants men"
$ws.Range("A36").Value = "mens black baseball pants"
$ws.Range("A37").Value = "knee pad protector"
$ws.Range("A38").Value = "compression pants boys basketball"
$ws.Range("A39").Value = "black leggings men"
$ws.Range("A40").Value = "athletic leggings men"
$ws.Range("A41").Value = "youth hex knee pads"
$ws.Range("A42").Value = "wrestling youth knee pads"
$ws.Range("A43").Value = "adult volleyball knee pads"
$ws.Range("A44").Value = "basketball shorts with pads"
$ws.Range("A45").Value = "mens wrestling tights"
$ws.Range("A46").Value = "softball sliding pants"
$ws.Range("A47").Value = "black baseball pants"
$ws.Range("A48").Value = "wrestling tights for men"
$ws.Range("A49").Value = "mens running capris"
$ws.Range("A50").Value = "mens football pants with pads"
$ws.Range("A51").Value = "soccer pants"
$ws.Range("A52").Value = "baseball sliding shorts boys"
$ws.Range("A53").Value = "compression knee guards"
$ws.Range("A54").Value = "knees protector"
$ws.Range("A55").Value = "baseball pants youth"
$ws.Range("A56").Value = "mens basketball shorts black"
$ws.Range("A57").Value = "mens compression"
$ws.Range("A58").Value = "wrestling clothes for men"
$ws.Range("A59").Value = "mens sports pants"
$ws.Range("A60").Value = "boys gym pants"
$ws.Range("A61").Value = "knee pads football adult"
$ws.Range("A62").Value = "leggings men short"
$ws.Range("A63").Value = "running pants men tall"
$ws.Range("A64").Value = "boys running tights youth"
$ws.Range("A65").Value = "softball pants for girls youth"
$ws.Range("A66").Value = "little boys athletic pants"
$ws.Range("A67").Value = "boys running pants size"
$ws.Range("A68").Value = "football youth pants"
$ws.Range("A69").Value = "comfortable knee pads"
$ws.Range("A70").Value = "knee sleeve with pad"
$ws.Range("A71").Value = "women compression tights"
$ws.Range("A72").Value = "knee pads toddler"
$ws.Range("A73").Value = "knee pads basketball kids"
$ws.Range("A74").Value = "knee protector for kids"
$ws.Range("A75").Value = "dark green knee pads for basketball"
$ws.Range("A76").Value = "skateboarding knee pads youth"
$ws.Range("A77").Value = "compression pants with padding basketball"
$ws.Range("A78").Value = "firefighter compression pants"
$ws.Range("A79").Value = "skateboard knee and elbow pads youth"
$ws.Range("A80").Value = "skateboard youth knee pads"
$ws.Range("A81").Value = "warm up pants men"
$ws.Range("A82").Value = "mc david knee pad"
$ws.Range("A83").Value = "compression knee pads men"
$ws.Range("A84").Value = "men basketball pants"
$ws.Range("A85").Value = "premium knee pad"
$ws.Range("A86").Value = "kids compression pants for basketball"
$ws.Range("A87").Value = "toddler knee pad"
$ws.Range("A88").Value = "padded knee sleeves men"
$ws.Range("A89").Value = "mtb knee pads for men"
$ws.Range("A90").Value = "compression knee sleeves pads"
$ws.Range("A91").Value = "basketball knee pads leggings"
$ws.Range("A92").Value = "soccer compression pants men"
$ws.Range("A93").Value = "mens small leggings"
$ws.Range("A94").Value = "compressions pants mens"
$ws.Range("A95").Value = "capri basketball leggings for boys"
$ws.Range("A96").Value = "youth small black baseball pants"
$ws.Range("A97").Value = "mens black leggings"
$ws.Range("A98").Value = "basket ball knee pads youth"
$ws.Range("A99").Value = "bjj leggings men"
$ws.Range("A100").Value = "mens gym tights"
